$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the last row (25) - the sheet shrinks from 25 to 24 rows,
#    matching the new dimension A1:C24.
# ---------------------------------------------------------------------
$ws.Rows.Item(25).Delete()

# ---------------------------------------------------------------------
# 2) Update cell contents so the remaining rows (10, 13-24) show the
#    new label/value pairing.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C10").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Types of wasteManagement and Legislation in Brazil"
$ws.Range("C14").Value = "Types of wasteManagement and Legislation in Brazil"

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").Value = "01/01/2020"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Concepts and definitions of solid waste; sampling, characterization and classification of solid waste; collection, conditioning, transportation, treatment and final disposal of solid waste; the National Policy on Solid Waste; model of solid waste management."
$ws.Range("C16").Value = "Concepts and definitions of solid waste; sampling, characterization and classification of solid waste; collection, conditioning, transportation, treatment and final disposal of solid waste; the National Policy on Solid Waste; model of solid waste management."

$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C18").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

$ws.Range("A19").Value = "Critério:"

$ws.Range("A20").Value = "Norma de recuperação:"

$ws.Range("A21").Value = "Bibliografia:"

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOB1202 -  Introdução ao Gerenciamento de Projetos Ambientais  (Requisito)`n"
$ws.Range("C23").Value = "LOB1202 -  Introdução ao Gerenciamento de Projetos Ambientais  (Requisito)`n"

$ws.Range("B24").Value = "LOB1232 -  Licenciamento Ambiental  (Requisito)`n"
$ws.Range("C24").Value = "LOB1232 -  Licenciamento Ambiental  (Requisito)`n"

# ---------------------------------------------------------------------
# 3) Fix up row heights so they match the new layout.
#    Rows 17 and 22 go back to the sheet's default (no custom height),
#    while rows 13, 15, 18, 21 and 23 need an explicit custom height.
# ---------------------------------------------------------------------
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()

$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
